# Fixed the units again - data is normalized to different (and different-per-column)
# powers of ten: column C gets an extra 10^(-3) factor, column D an extra 10^(-4) factor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 keeps its own (non-shared) formulas.
$ws.Range("C2").Formula = "=A2/1.4765679173556 * 10^(-3)"
$ws.Range("D2").Formula = "=B2/1.4765679173556 * 10^(-4)"

# Rows 3-66 become one shared-formula block per column (C and D can no longer
# share a single formula group because the multiplier now differs per column).
$ws.Range("C3:C66").Formula = "=A3/1.4765679173556 * 10^(-3)"
$ws.Range("D3:D66").Formula = "=B3/1.4765679173556 * 10^(-4)"

# Rows 67-68 form their own trailing shared-formula block, same reasoning.
$ws.Range("C67:C68").Formula = "=A67/1.4765679173556 * 10^(-3)"
$ws.Range("D67:D68").Formula = "=B67/1.4765679173556 * 10^(-4)"

# Column widths for C:D were set to a fixed best-fit-style width.
$ws.Range("C:D").ColumnWidth = 11.1

# Selection moved from the single cell D68 to the whole D2:D68 range (active cell D2).
$ws.Range("D2:D68").Select()
